$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.512.83"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").Value = "2.983.28"
$ws.Range("E3").Value = "  +2.75%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.99%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.545"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.69%  "

$ws.Range("E11").Value = "  +0.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0844"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.74%  "

$ws.Range("D14").Value = "3.455.61"
$ws.Range("E14").Value = "  +2.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.60%  "

$ws.Range("D16").Value = "2.980.70"
$ws.Range("E16").Value = "  +2.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.971"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.38%  "

$ws.Range("D18").Value = "51.594.82"
$ws.Range("E18").Value = "  +1.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  +2.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.32%  "

$ws.Range("E25").Value = "  +4.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.171"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +19.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.31%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.66%  "

$ws.Range("E31").Value = "  +5.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.77%  "

$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.08%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.44%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0461"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.52%  "

$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.69%  "

$ws.Range("E40").Value = "  -5.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "123.56"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.282"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +21.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.00%  "

$ws.Range("D48").Value = "2.052.92"
$ws.Range("E48").Value = "  -0.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0355"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.70%  "

$ws.Range("E51").Value = "  +4.00%  "
